# Update: CODECOP 4.0.0 report refresh — reprocessed run with updated timings,
# complexity threshold (5 -> 4), newly-surfaced complex PL/SQL units, and
# refreshed per-file "issues density" (column R) on the Files sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Defined name UnitsTable now spans the populated PLSQLUnits rows (1:8).
# ---------------------------------------------------------------------------
$unitsName = $wb.Names.Item("UnitsTable")
$unitsName.RefersTo = "=PLSQLUnits!`$A`$1:`$H`$8"

# ---------------------------------------------------------------------------
# 2. Summary sheet: complexity parameter + processing timestamps/duration.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

# "complexity" parameter dropped from 5 to 4 (stored as text, same as before).
$summary.Range("B21").NumberFormat = "@"
$summary.Range("B21").Value2 = "4"

# Start / end of processing (Excel date-time serials) and elapsed seconds.
$summary.Range("B36").Value2 = 44287.835625
$summary.Range("B37").Value2 = 44287.83577546296
$summary.Range("B38").Value2 = 10.031999999999998

# ---------------------------------------------------------------------------
# 3. PLSQLUnits sheet: newly reported complex PL/SQL units (rows 2-8).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PLSQLUnits")

$newRows = @(
    @('guidelines/guideline_1040_04.sql', 'AnonymousPlsqlBlock', 10.0, 36.0, 1.0, 6.0, 29.0, 14.0, 5.0, 349.0333754971396, 94.57346819440446),
    @('guidelines/guideline_4370_45.sql', 'AnonymousPlsqlBlock', 14.0, 30.0, 0.0, 3.0, 27.0, 13.0, 5.0, 411.1982937621106, 83.45140948052318),
    @('guidelines/guideline_4310_39.sql', 'my_package.password_check', 26.0, 24.0, 0.0, 3.0, 21.0, 10.0, 5.0, 491.54240635418904, 86.13827722867563),
    @('guidelines/guideline_4310_39.sql', 'my_package.password_check', 63.0, 24.0, 0.0, 3.0, 21.0, 10.0, 5.0, 491.54240635418904, 86.13827722867563),
    @('guidelines/guideline_4320_40.sql', 'AnonymousPlsqlBlock', 12.0, 21.0, 0.0, 3.0, 18.0, 9.0, 5.0, 288.85263754543286, 91.06596991130587),
    @('guidelines/guideline_4320_40.sql', 'AnonymousPlsqlBlock', 41.0, 26.0, 0.0, 3.0, 23.0, 9.0, 5.0, 346.1295543881475, 86.66540775801556),
    @('guidelines/guideline_4370_45.sql', 'AnonymousPlsqlBlock', 52.0, 26.0, 0.0, 3.0, 23.0, 9.0, 5.0, 346.1295543881475, 86.66540775801556)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $i + 2
    $rowValues = $newRows[$i]
    $ws.Range("A$r").Value2 = $rowValues[0]
    $ws.Range("B$r").Value2 = $rowValues[1]
    $ws.Range("C$r").Value2 = $rowValues[2]
    $ws.Range("D$r").Value2 = $rowValues[3]
    $ws.Range("E$r").Value2 = $rowValues[4]
    $ws.Range("F$r").Value2 = $rowValues[5]
    $ws.Range("G$r").Value2 = $rowValues[6]
    $ws.Range("H$r").Value2 = $rowValues[7]
    $ws.Range("I$r:K$r").NumberFormat = "0"
    $ws.Range("I$r").Value2 = $rowValues[8]
    $ws.Range("J$r").Value2 = $rowValues[9]
    $ws.Range("K$r").Value2 = $rowValues[10]
}

# ---------------------------------------------------------------------------
# 4. Files sheet: refreshed "issues density" (column R) for re-scanned files.
# ---------------------------------------------------------------------------
$files = $wb.Worksheets.Item("Files")

$rColumnUpdates = @{
    2 = 0.031
    4 = 0.034
    5 = 0.05
    7 = 0.073
    8 = 0.036
    9 = 0.053
    10 = 0.057
    11 = 0.018
    12 = 0.021
    13 = 0.02
    15 = 0.02
    16 = 0.016
    17 = 0.035
    18 = 0.03
    19 = 0.052
    20 = 0.021
    21 = 0.025
    22 = 0.013
    23 = 0.083
    24 = 0.021
    25 = 0.021
    26 = 0.018
    28 = 0.03
    29 = 0.041
    30 = 0.035
    31 = 0.019
    32 = 0.025
    33 = 0.028
    34 = 0.018
    35 = 0.026
    36 = 0.018
    37 = 0.037
    39 = 0.026
    40 = 0.015
    41 = 0.019
    43 = 0.02
    44 = 0.019
    45 = 0.02
    46 = 0.024
    47 = 0.029
    48 = 0.017
    49 = 0.019
    50 = 0.015
    53 = 0.022
    54 = 1.846
    55 = 0.03
    56 = 0.097
    58 = 0.012
    59 = 0.021
    60 = 0.031
    61 = 0.016
    62 = 0.029
    63 = 0.027
    64 = 0.03
    65 = 0.016
    67 = 0.03
    68 = 0.022
    70 = 0.031
    71 = 0.044
    72 = 0.049
    73 = 0.036
    74 = 0.036
    75 = 0.028
    76 = 0.017
    77 = 0.019
    78 = 0.025
    79 = 0.037
    80 = 0.022
    82 = 0.013
    83 = 0.014
    84 = 0.03
    85 = 0.051
    86 = 0.038
    87 = 0.035
    88 = 0.034
    89 = 0.052
    90 = 0.021
    91 = 0.016
    92 = 0.016
    93 = 0.016
    94 = 0.023
    95 = 0.026
    96 = 0.03
    98 = 0.022
    99 = 0.023
    100 = 0.021
    101 = 0.017
    102 = 0.02
    103 = 0.016
    104 = 0.015
    105 = 5.227
    106 = 0.014
    107 = 0.017
    109 = 0.013
    111 = 0.027
    112 = 0.02
    113 = 0.016
    115 = 0.032
}

foreach ($row in $rColumnUpdates.Keys) {
    $files.Range("R$row").Value2 = $rColumnUpdates[$row]
}
